$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich text cells): bump the issue number and the two
# report-week dates. We only touch the exact characters that changed so the
# surrounding run layout / wording stays intact.
# ---------------------------------------------------------------------------

# A8: "Volume 32   Number  3" -> "Volume 32   Number  4"
$ws.Range("A8").Characters(21, 1).Text = "4"

# C9: "Report Covering the Week  1/13/2025  Through  1/19/2025"
#  -> "Report Covering the Week  1/20/2025  Through  1/26/2025"
$ws.Range("C9").Characters(27, 9).Text = "1/20/2025"
$ws.Range("C9").Characters(47, 9).Text = "1/26/2025"

# ---------------------------------------------------------------------------
# Weekly crime-stat grid updates (rows 16-28). Most cells simply get a new
# number; a handful of cells change from a numeric value to the "no data"
# text placeholders ("0" / "***.*") or vice versa, which also needs the
# cell's style switched between the text style (13) and the matching numeric
# style (14 integer / 15 percent) so the workbook matches the source types.
# ---------------------------------------------------------------------------

function Set-NumberCell($addr, $value, $numericFormatSource) {
    # Used when a cell that previously held placeholder text now holds a
    # real number: write the value, then copy over the number format/style
    # from a cell that already has the right numeric style.
    $ws.Range($addr).Value = $value
    $ws.Range($numericFormatSource).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

function Set-TextCell($addr, $text, $textFormatSource) {
    # Used when a cell that previously held a real number now holds one of
    # the placeholder text strings: paste the whole reference cell in first
    # (so the cell becomes a genuine text cell), then reapply its format so
    # the cell ends up on the same shared text style as its neighbours.
    $ws.Range($textFormatSource).Copy() | Out-Null
    $ws.Paste($ws.Range($addr)) | Out-Null
    $ws.Range($textFormatSource).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# --- Row 16 ---
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 9
$ws.Range("K16").Value = -44.444444444444
$ws.Range("L16").Value = -61.538461538461
$ws.Range("M16").Value = -72.222222222222
$ws.Range("N16").Value = -92.957746478873

# --- Row 17 ---
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = -12.5
$ws.Range("I17").Value = 5
$ws.Range("J17").Value = 8
$ws.Range("K17").Value = -37.5
$ws.Range("L17").Value = -37.5
$ws.Range("M17").Value = -16.666666666666
$ws.Range("N17").Value = -66.666666666666

# --- Row 18 ---
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -25
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 19
$ws.Range("K18").Value = -42.105263157894
$ws.Range("L18").Value = -62.068965517241
$ws.Range("M18").Value = -52.173913043478
$ws.Range("N18").Value = -80

# --- Row 19 ---
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -29.411764705882
$ws.Range("F19").Value = 70
$ws.Range("H19").Value = -11.392405063291
$ws.Range("I19").Value = 63
$ws.Range("J19").Value = 70
$ws.Range("K19").Value = -10
$ws.Range("L19").Value = -37
$ws.Range("M19").Value = -10
$ws.Range("N19").Value = -65.384615384615

# --- Row 20 --- (D20, E20 go from placeholder text to real numbers)
Set-NumberCell "D20" 1 "F20"
Set-NumberCell "E20" -100 "H20"
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -75
$ws.Range("J20").Value = 3
$ws.Range("K20").Value = -66.666666666666
$ws.Range("N20").Value = -98.461538461538

# --- Row 21 (TOTAL row) ---
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = -37.037037037037
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 122
$ws.Range("H21").Value = -18.852459016393
$ws.Range("I21").Value = 85
$ws.Range("J21").Value = 109
$ws.Range("K21").Value = -22.018348623853
$ws.Range("L21").Value = -44.078947368421
$ws.Range("M21").Value = -27.966101694915
$ws.Range("N21").Value = -78.149100257069

# --- Row 22 --- (D22, E22 go from real numbers to placeholder text)
Set-TextCell "D22" "0" "C20"
Set-TextCell "E22" "***.*" "E14"
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -20
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = -20
$ws.Range("L22").Value = 100

# --- Row 24 ---
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -28.571428571428
$ws.Range("F24").Value = 112
$ws.Range("G24").Value = 104
$ws.Range("H24").Value = 7.692307692307
$ws.Range("I24").Value = 101
$ws.Range("J24").Value = 92
$ws.Range("K24").Value = 9.782608695652
$ws.Range("L24").Value = -20.472440944881
$ws.Range("M24").Value = 7.446808510638

# --- Row 25 ---
$ws.Range("C25").Value = 17
$ws.Range("E25").Value = -32
$ws.Range("F25").Value = 78
$ws.Range("H25").Value = -11.363636363636
$ws.Range("I25").Value = 70
$ws.Range("J25").Value = 76
$ws.Range("K25").Value = -7.894736842105
$ws.Range("L25").Value = -28.571428571428

# --- Row 26 ---
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 133.333333333333
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = -11.111111111111
$ws.Range("I26").Value = 21
$ws.Range("J26").Value = 24
$ws.Range("K26").Value = -12.5
$ws.Range("M26").Value = 133.333333333333

# --- Row 28 --- (C28 goes from a real number to placeholder text)
Set-TextCell "C28" "0" "D28"
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 66.666666666666
$ws.Range("I28").Value = 5
$ws.Range("K28").Value = 150
$ws.Range("L28").Value = -37.5
